$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Product table rows 12-80: Qty (A), Barcode (B), Product Name (C)
# Rows 1-11 are unchanged by this edit.
$data = @(
    ,@(12, 3, '782796018734', 'Fruit Cheese & Almond Bites')
    ,@(13, 3, '782796018598', 'Fruit & Cheese Bites')
    ,@(14, 1, '854426008139', 'Super Fruit Medley')
    ,@(15, 1, '854426008252', 'Berry Blend')
    ,@(16, 1, '854426008245', 'Melon Trio')
    ,@(17, 1, '854426008436', 'Watermelon Spears')
    ,@(18, 1, '854426008429', 'Pineapple Spears')
    ,@(19, 2, '854426008337', 'Sliced Apples')
    ,@(20, 1, '000000598392', 'GNG White Half Turkey Bacon Club Sub')
    ,@(21, 1, '000000523899', 'GNG White Half Turkey Sub')
    ,@(22, 2, '000000520690', 'GNG White Whole Turkey Sub')
    ,@(23, 1, '000000523851', 'GNG Wheat Half Ham Sub')
    ,@(24, 2, '854426008146', 'Red Grape Cup')
    ,@(25, 2, '854426008122', 'Seasonal Fruit Bowl')
    ,@(26, 1, '030223036135', 'Chicken Caesar Salad')
    ,@(27, 1, '030223071174', 'Chicken Cobb Salad')
    ,@(28, 2, '000000523882', 'GNG White Half Italian Sub')
    ,@(29, 2, '000000520706', 'GNG White Whole Italian Sub')
    ,@(30, 1, '000000520713', 'GNG Wheat Whole Stack Sub')
    ,@(31, 2, '000000541138', 'GNG Southwest Chicken Wrap')
    ,@(32, 1, '000000532099', 'GNG Turkey Wrap')
    ,@(33, 1, '030223036135', 'Chicken Caesar Salad')
    ,@(34, 1, '030223071174', 'Chicken Cobb Salad')
    ,@(35, 1, '030223071181', 'Chef Salad Ham/Turkey')
    ,@(36, 1, '040822011907', 'Classic Hummus w/Pretzels')
    ,@(37, 1, '040822011990', 'Red Pepper Hummus w/Pretzels')
    ,@(38, 1, '046675026976', 'Yo Crunch Vanilla M&M')
    ,@(39, 1, '046675027021', 'Yo Crunch Vanilla Oreo')
    ,@(40, 1, '894700010052', 'Chobani Greek Blueberry')
    ,@(41, 1, '889470001003', 'Chobani Greek Strawberry')
    ,@(42, 1, '854426008078', 'Strawberry & Banana Yogurt Parfait')
    ,@(43, 1, '782796018543', 'Strawberry Yogurt Parfait')
    ,@(44, 1, '040697640806', 'Triple Chocolate Bunt Cake')
    ,@(45, 1, '040697719097', 'Carrot Cake w/Cream Cheese')
    ,@(46, 2, '786162411716', 'Smart Water Alkaline 12oz Can')
    ,@(47, 2, '786162411709', 'Smart Water 12oz Can')
    ,@(48, 2, '026400700043', 'Darigold 1% Chocolate MILK 8oz')
    ,@(49, 2, '026400700098', 'Darigold 2% Milk 8oz')
    ,@(50, 1, '851554006089', 'Noka Strawberry Pineapple')
    ,@(51, 1, '851554006454', 'Noka Superfood')
    ,@(52, 1, '074329123143', 'Oh Snap Cranberry')
    ,@(53, 1, '074329123256', 'Oh Snap Sassy Bites')
    ,@(54, 1, '074329123382', 'Oh Snap Dilly Bites')
    ,@(55, 1, '074329123393', 'Oh Snap Hottie Bites')
    ,@(56, 2, '049000400441', 'Coke 20oz')
    ,@(57, 2, '049000400458', 'Coke Diet 20oz')
    ,@(58, 2, '049000407648', 'Sprite 20oz')
    ,@(59, 2, '078000082401', 'Dr Pepper 20oz')
    ,@(60, 2, '012000001291', 'Pepsi 20oz')
    ,@(61, 2, '012000001314', 'Mountain Dew 20oz')
    ,@(62, 1, '898999012698', 'Vita Coco w/Pulp 16oz Can')
    ,@(63, 1, '898999012704', 'Vita Coco w/Mango 16oz Can')
    ,@(64, 2, '853004004020', 'Core Water 20oz')
    ,@(65, 2, '898999010229', 'Vita Coco 17oz Btl')
    ,@(66, 2, '049000409772', 'Dasani 20oz')
    ,@(67, 1, '049000407907', 'Powerade Mountain Blast 20oz')
    ,@(68, 1, '049000403718', 'Powerade Fruit Punch 20oz')
    ,@(69, 2, '025000062193', 'Minute Maid Orange Juice')
    ,@(70, 2, '025000061523', 'Minute Maid Apple Juice')
    ,@(71, 2, '811620022002', 'Core Power Chocolate')
    ,@(72, 2, '811620022033', 'Core Power Strawberry')
    ,@(73, 1, '811620022019', 'Core Power Vanilla')
    ,@(74, 2, '049000172386', 'Dunkin Donuts Mocha 13.7oz')
    ,@(75, 1, '049000172393', 'Dunkin Donuts Vanilla 13.7oz')
    ,@(76, 1, '710779006539', 'Lean Body Chocolate Peanutbutter 14oz')
    ,@(77, 1, '710779770515', 'Lean Body Strawberry 14oz')
    ,@(78, 2, '611269917475', 'Red Bull 12oz')
    ,@(79, 1, '611269002072', 'Red Bull SF Watermelon 12oz')
    ,@(80, 1, '611269002157', 'Red Bull SF Strawberry Apricot 12oz')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Rows 78-79 use an explicit black font color (distinct style) in the source file;
# row 80 uses the normal/default style.
$ws.Range("B78:C79").Font.Color = 0

# Re-fit the A and C columns since content length changed
$ws.Columns("A:C").AutoFit() | Out-Null

# Restore view state (scroll position + active selection) to match the saved workbook
$ws.Range("B78").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 53

